$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 73.38544233333333
$ws.Range("H2").Value = 220.156327
$ws.Range("I2").Value = 0.1214979676060253
$ws.Range("J2").Value = 0.1214979676060253
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6187893333333333
$ws.Range("N2").Value = 1.856368
$ws.Range("O2").Value = 0.1765034355725207
$ws.Range("P2").Value = 0.1765034355725208
$ws.Range("Q2").Value = 45.41012893781511
$ws.Range("R2").Value = 408.691160440336
$ws.Range("S2").Value = 0.0214448086975423
$ws.Range("T2").Value = 0.0214448086975423
$ws.Range("G3").Value = 73.38544233333333
$ws.Range("H3").Value = 220.156327
$ws.Range("I3").Value = 0.1214979676060253
$ws.Range("J3").Value = 0.1214979676060253
$ws.Range("N3").Value = 6.701951000000001
$ws.Range("O3").Value = 0.6372213788099619
$ws.Range("P3").Value = 0.6372213788099619
$ws.Range("Q3").Value = 163.9418795437752
$ws.Range("R3").Value = 1475.476915893977
$ws.Range("S3").Value = 0.07742110244051953
$ws.Range("T3").Value = 0.07742110244051952
$ws.Range("G4").Value = 73.38544233333333
$ws.Range("H4").Value = 220.156327
$ws.Range("I4").Value = 0.1214979676060253
$ws.Range("J4").Value = 0.1214979676060253
$ws.Range("M4").Value = 0.207158
$ws.Range("N4").Value = 0.6214740000000001
$ws.Range("O4").Value = 0.05908973658186135
$ws.Range("P4").Value = 0.05908973658186135
$ws.Range("Q4").Value = 15.20238146288867
$ws.Range("R4").Value = 136.821433165998
$ws.Range("S4").Value = 0.007179282901071558
$ws.Range("T4").Value = 0.007179282901071557
$ws.Range("G5").Value = 73.38544233333333
$ws.Range("H5").Value = 220.156327
$ws.Range("I5").Value = 0.1214979676060253
$ws.Range("J5").Value = 0.1214979676060253
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1889926666666667
$ws.Range("N5").Value = 0.566978
$ws.Range("O5").Value = 0.05390825789608347
$ws.Range("P5").Value = 0.05390825789608347
$ws.Range("Q5").Value = 13.86931044108955
$ws.Range("R5").Value = 124.823793969806
$ws.Range("S5").Value = 0.006549743771555607
$ws.Range("T5").Value = 0.006549743771555605
$ws.Range("G6").Value = 73.38544233333333
$ws.Range("H6").Value = 220.156327
$ws.Range("I6").Value = 0.1214979676060253
$ws.Range("J6").Value = 0.1214979676060253
$ws.Range("M6").Value = 0.2568966666666667
$ws.Range("N6").Value = 0.77069
$ws.Range("O6").Value = 0.07327719113957255
$ws.Range("P6").Value = 0.07327719113957255
$ws.Range("Q6").Value = 18.85247551729222
$ws.Range("R6").Value = 169.67227965563
$ws.Range("S6").Value = 0.008903029795336311
$ws.Range("T6").Value = 0.008903029795336309
$ws.Range("I7").Value = 0.3924995450689984
$ws.Range("J7").Value = 0.3924995450689983
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6187893333333333
$ws.Range("N7").Value = 1.856368
$ws.Range("O7").Value = 0.1765034355725207
$ws.Range("P7").Value = 0.1765034355725208
$ws.Range("Q7").Value = 146.6975563526471
$ws.Range("R7").Value = 1320.278007173824
$ws.Range("S7").Value = 0.06927751816532966
$ws.Range("T7").Value = 0.06927751816532965
$ws.Range("I8").Value = 0.3924995450689984
$ws.Range("J8").Value = 0.3924995450689983
$ws.Range("N8").Value = 6.701951000000001
$ws.Range("O8").Value = 0.6372213788099619
$ws.Range("P8").Value = 0.6372213788099619
$ws.Range("Q8").Value = 529.6147285964744
$ws.Range("R8").Value = 4766.532557368269
$ws.Range("S8").Value = 0.2501091012911499
$ws.Range("T8").Value = 0.2501091012911499
$ws.Range("I9").Value = 0.3924995450689984
$ws.Range("J9").Value = 0.3924995450689983
$ws.Range("M9").Value = 0.207158
$ws.Range("N9").Value = 0.6214740000000001
$ws.Range("O9").Value = 0.05908973658186135
$ws.Range("P9").Value = 0.05908973658186135
$ws.Range("Q9").Value = 49.11133845051468
$ws.Range("R9").Value = 442.0020460546321
$ws.Range("S9").Value = 0.02319269472662753
$ws.Range("T9").Value = 0.02319269472662753
$ws.Range("I10").Value = 0.3924995450689984
$ws.Range("J10").Value = 0.3924995450689983
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1889926666666667
$ws.Range("N10").Value = 0.566978
$ws.Range("O10").Value = 0.05390825789608347
$ws.Range("P10").Value = 0.05390825789608347
$ws.Range("Q10").Value = 44.80484855681156
$ws.Range("R10").Value = 403.243637011304
$ws.Range("S10").Value = 0.02115896669967501
$ws.Range("T10").Value = 0.021158966699675
$ws.Range("I11").Value = 0.3924995450689984
$ws.Range("J11").Value = 0.3924995450689983
$ws.Range("M11").Value = 0.2568966666666667
$ws.Range("N11").Value = 0.77069
$ws.Range("O11").Value = 0.07327719113957255
$ws.Range("P11").Value = 0.07327719113957255
$ws.Range("Q11").Value = 60.90297813010223
$ws.Range("R11").Value = 548.1268031709201
$ws.Range("S11").Value = 0.02876126418621627
$ws.Range("T11").Value = 0.02876126418621626
$ws.Range("G12").Value = 138.1628113333333
$ws.Range("H12").Value = 414.488434
$ws.Range("I12").Value = 0.2287442882675098
$ws.Range("J12").Value = 0.2287442882675098
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.6187893333333333
$ws.Range("N12").Value = 1.856368
$ws.Range("O12").Value = 0.1765034355725207
$ws.Range("P12").Value = 0.1765034355725208
$ws.Range("Q12").Value = 85.49367391641245
$ws.Range("R12").Value = 769.443065247712
$ws.Range("S12").Value = 0.04037415274680654
$ws.Range("T12").Value = 0.04037415274680653
$ws.Range("G13").Value = 138.1628113333333
$ws.Range("H13").Value = 414.488434
$ws.Range("I13").Value = 0.2287442882675098
$ws.Range("J13").Value = 0.2287442882675098
$ws.Range("N13").Value = 6.701951000000001
$ws.Range("O13").Value = 0.6372213788099619
$ws.Range("P13").Value = 0.6372213788099619
$ws.Range("Q13").Value = 308.653463859415
$ws.Range("R13").Value = 2777.881174734734
$ws.Range("S13").Value = 0.145760750764726
$ws.Range("T13").Value = 0.145760750764726
$ws.Range("G14").Value = 138.1628113333333
$ws.Range("H14").Value = 414.488434
$ws.Range("I14").Value = 0.2287442882675098
$ws.Range("J14").Value = 0.2287442882675098
$ws.Range("M14").Value = 0.207158
$ws.Range("N14").Value = 0.6214740000000001
$ws.Range("O14").Value = 0.05908973658186135
$ws.Range("P14").Value = 0.05908973658186135
$ws.Range("Q14").Value = 28.62153167019067
$ws.Range("R14").Value = 257.593785031716
$ws.Range("S14").Value = 0.01351643973833251
$ws.Range("T14").Value = 0.01351643973833251
$ws.Range("G15").Value = 138.1628113333333
$ws.Range("H15").Value = 414.488434
$ws.Range("I15").Value = 0.2287442882675098
$ws.Range("J15").Value = 0.2287442882675098
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.1889926666666667
$ws.Range("N15").Value = 0.566978
$ws.Range("O15").Value = 0.05390825789608347
$ws.Range("P15").Value = 0.05390825789608347
$ws.Range("Q15").Value = 26.11175814805022
$ws.Range("R15").Value = 235.005823332452
$ws.Range("S15").Value = 0.01233120608418098
$ws.Range("T15").Value = 0.01233120608418098
$ws.Range("G16").Value = 138.1628113333333
$ws.Range("H16").Value = 414.488434
$ws.Range("I16").Value = 0.2287442882675098
$ws.Range("J16").Value = 0.2287442882675098
$ws.Range("M16").Value = 0.2568966666666667
$ws.Range("N16").Value = 0.77069
$ws.Range("O16").Value = 0.07327719113957255
$ws.Range("P16").Value = 0.07327719113957255
$ws.Range("Q16").Value = 35.49356568882889
$ws.Range("R16").Value = 319.44209119946
$ws.Range("S16").Value = 0.0167617389334638
$ws.Range("T16").Value = 0.0167617389334638
$ws.Range("G17").Value = 49.051656
$ws.Range("H17").Value = 147.154968
$ws.Range("I17").Value = 0.08121060965524597
$ws.Range("J17").Value = 0.08121060965524596
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.6187893333333333
$ws.Range("N17").Value = 1.856368
$ws.Range("O17").Value = 0.1765034355725207
$ws.Range("P17").Value = 0.1765034355725208
$ws.Range("Q17").Value = 30.352641515136
$ws.Range("R17").Value = 273.173773636224
$ws.Range("S17").Value = 0.01433395160908984
$ws.Range("T17").Value = 0.01433395160908984
$ws.Range("G18").Value = 49.051656
$ws.Range("H18").Value = 147.154968
$ws.Range("I18").Value = 0.08121060965524597
$ws.Range("J18").Value = 0.08121060965524596
$ws.Range("N18").Value = 6.701951000000001
$ws.Range("O18").Value = 0.6372213788099619
$ws.Range("P18").Value = 0.6372213788099619
$ws.Range("Q18").Value = 109.580598326952
$ws.Range("R18").Value = 986.2253849425681
$ws.Range("S18").Value = 0.05174913665851344
$ws.Range("T18").Value = 0.05174913665851343
$ws.Range("G19").Value = 49.051656
$ws.Range("H19").Value = 147.154968
$ws.Range("I19").Value = 0.08121060965524597
$ws.Range("J19").Value = 0.08121060965524596
$ws.Range("M19").Value = 0.207158
$ws.Range("N19").Value = 0.6214740000000001
$ws.Range("O19").Value = 0.05908973658186135
$ws.Range("P19").Value = 0.05908973658186135
$ws.Range("Q19").Value = 10.161442953648
$ws.Range("R19").Value = 91.45298658283201
$ws.Range("S19").Value = 0.00479871353218085
$ws.Range("T19").Value = 0.004798713532180849
$ws.Range("G20").Value = 49.051656
$ws.Range("H20").Value = 147.154968
$ws.Range("I20").Value = 0.08121060965524597
$ws.Range("J20").Value = 0.08121060965524596
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.1889926666666667
$ws.Range("N20").Value = 0.566978
$ws.Range("O20").Value = 0.05390825789608347
$ws.Range("P20").Value = 0.05390825789608347
$ws.Range("Q20").Value = 9.270403271856001
$ws.Range("R20").Value = 83.43362944670399
$ws.Range("S20").Value = 0.004377922489193166
$ws.Range("T20").Value = 0.004377922489193165
$ws.Range("G21").Value = 49.051656
$ws.Range("H21").Value = 147.154968
$ws.Range("I21").Value = 0.08121060965524597
$ws.Range("J21").Value = 0.08121060965524596
$ws.Range("M21").Value = 0.2568966666666667
$ws.Range("N21").Value = 0.77069
$ws.Range("O21").Value = 0.07327719113957255
$ws.Range("P21").Value = 0.07327719113957255
$ws.Range("Q21").Value = 12.60120692088
$ws.Range("R21").Value = 113.41086228792
$ws.Range("S21").Value = 0.005950885366268675
$ws.Range("T21").Value = 0.005950885366268675
$ws.Range("G22").Value = 106.3337146666667
$ws.Range("H22").Value = 319.001144
$ws.Range("I22").Value = 0.1760475894022206
$ws.Range("J22").Value = 0.1760475894022206
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.6187893333333333
$ws.Range("N22").Value = 1.856368
$ws.Range("O22").Value = 0.1765034355725207
$ws.Range("P22").Value = 0.1765034355725208
$ws.Range("Q22").Value = 65.79816840944355
$ws.Range("R22").Value = 592.183515684992
$ws.Range("S22").Value = 0.03107300435375243
$ws.Range("T22").Value = 0.03107300435375243
$ws.Range("G23").Value = 106.3337146666667
$ws.Range("H23").Value = 319.001144
$ws.Range("I23").Value = 0.1760475894022206
$ws.Range("J23").Value = 0.1760475894022206
$ws.Range("N23").Value = 6.701951000000001
$ws.Range("O23").Value = 0.6372213788099619
$ws.Range("P23").Value = 0.6372213788099619
$ws.Range("Q23").Value = 237.5477817813272
$ws.Range("R23").Value = 2137.930036031944
$ws.Range("S23").Value = 0.112181287655053
$ws.Range("T23").Value = 0.112181287655053
$ws.Range("G24").Value = 106.3337146666667
$ws.Range("H24").Value = 319.001144
$ws.Range("I24").Value = 0.1760475894022206
$ws.Range("J24").Value = 0.1760475894022206
$ws.Range("M24").Value = 0.207158
$ws.Range("N24").Value = 0.6214740000000001
$ws.Range("O24").Value = 0.05908973658186135
$ws.Range("P24").Value = 0.05908973658186135
$ws.Range("Q24").Value = 22.02787966291734
$ws.Range("R24").Value = 198.250916966256
$ws.Range("S24").Value = 0.0104026056836489
$ws.Range("T24").Value = 0.0104026056836489
$ws.Range("G25").Value = 106.3337146666667
$ws.Range("H25").Value = 319.001144
$ws.Range("I25").Value = 0.1760475894022206
$ws.Range("J25").Value = 0.1760475894022206
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 0.1889926666666667
$ws.Range("N25").Value = 0.566978
$ws.Range("O25").Value = 0.05390825789608347
$ws.Range("P25").Value = 0.05390825789608347
$ws.Range("Q25").Value = 20.09629229142578
$ws.Range("R25").Value = 180.866630622832
$ws.Range("S25").Value = 0.00949041885147872
$ws.Range("T25").Value = 0.009490418851478716
$ws.Range("G26").Value = 106.3337146666667
$ws.Range("H26").Value = 319.001144
$ws.Range("I26").Value = 0.1760475894022206
$ws.Range("J26").Value = 0.1760475894022206
$ws.Range("M26").Value = 0.2568966666666667
$ws.Range("N26").Value = 0.77069
$ws.Range("O26").Value = 0.07327719113957255
$ws.Range("P26").Value = 0.07327719113957255
$ws.Range("Q26").Value = 27.31677685215111
$ws.Range("R26").Value = 245.85099166936
$ws.Range("S26").Value = 0.01290027285828751
$ws.Range("T26").Value = 0.0129002728582875

Write-Output "Updated 305 cells"
